# New crime data collected — roll the CompStat "7th Precinct" weekly report
# forward by one week (Volume/Number + report-date header) and refresh the
# Week-to-Date / 28-Day / Year-to-Date / 2-Year crime figures and their
# computed percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 30   Number  18" -> "...19", and the report-week dates
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# ---------------------------------------------------------------------
# Helper constants for PasteSpecial (avoids relying on named enum members)
#   -4122 = xlPasteFormats
#   -4163 = xlPasteValues
# ---------------------------------------------------------------------
$xlPasteFormats = -4122
$xlPasteValues = -4163

# ---------------------------------------------------------------------
# Cells that flip between the numeric "#,##0" style and the text-styled
# placeholder ("0" / "***.*" shared strings used elsewhere in this sheet
# for zero-count / not-applicable cells). Copy the number format from a
# same-styled neighbour first so the underlying style index matches, then
# write the value.
# ---------------------------------------------------------------------

# Rape row (15): Week-to-Date 2023 count 1 -> no reported complaints ("0")
$ws.Range("D15").Copy()
$ws.Range("F15").PasteSpecial($xlPasteFormats)
$ws.Range("D15").Copy()
$ws.Range("F15").PasteSpecial($xlPasteValues)

# Transit (22): Week-to-Date 2023 count "0" -> 1 (now numeric)
$ws.Range("G14").Copy()
$ws.Range("C22").PasteSpecial($xlPasteFormats)
$ws.Range("C22").Value = 1
$null = 0 # (placeholder kept intentionally blank)

# Housing (30): Week-to-Date 2023 count 1 -> no reported complaints ("0")
$ws.Range("D15").Copy()
$ws.Range("F30").PasteSpecial($xlPasteFormats)
$ws.Range("D15").Copy()
$ws.Range("F30").PasteSpecial($xlPasteValues)

# G.L.A. (20): Week-to-Date 2023 "0" -> 1 (now numeric)
$ws.Range("G14").Copy()
$ws.Range("C20").PasteSpecial($xlPasteFormats)
$ws.Range("C20").Value = 1

# Transit (22): Week-to-Date 2022 "0" -> 1 (now numeric)
$ws.Range("G14").Copy()
$ws.Range("D22").PasteSpecial($xlPasteFormats)
$ws.Range("D22").Value = 1

# Transit (22): Week-to-Date %Chg "***.*" -> -100 (now numeric)
$ws.Range("H14").Copy()
$ws.Range("E22").PasteSpecial($xlPasteFormats)
$ws.Range("E22").Value = -100

# UCR Rape* (26): Week-to-Date 2023/2022 "0" -> 1 / 1, %Chg "***.*" -> 0
$ws.Range("G14").Copy()
$ws.Range("C26").PasteSpecial($xlPasteFormats)
$ws.Range("C26").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D26").PasteSpecial($xlPasteFormats)
$ws.Range("D26").Value = 1
$ws.Range("H14").Copy()
$ws.Range("E26").PasteSpecial($xlPasteFormats)
$ws.Range("E26").Value = 0

# UCR Rape* (26): 28-Day 2022 "0" -> 1, %Chg "***.*" -> 0
$ws.Range("G14").Copy()
$ws.Range("G26").PasteSpecial($xlPasteFormats)
$ws.Range("G26").Value = 1
$ws.Range("H14").Copy()
$ws.Range("H26").PasteSpecial($xlPasteFormats)
$ws.Range("H26").Value = 0

# Other Sex Crimes (27): Week-to-Date 2022 "0" -> 1, %Chg "***.*" -> 0
$ws.Range("G14").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$ws.Range("D27").Value = 1
$ws.Range("H14").Copy()
$ws.Range("E27").PasteSpecial($xlPasteFormats)
$ws.Range("E27").Value = 0

# ---------------------------------------------------------------------
# Plain numeric refreshes (same style/number-format, new counts & % chg)
# ---------------------------------------------------------------------

# Robbery (row 16)
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -25.423728813559
$ws.Range("L16").Value = 144.444444444444
$ws.Range("M16").Value = -6.382978723404
$ws.Range("N16").Value = -84.115523465704

# Fel. Assault (row 17)
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -44.444444444444
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 70
$ws.Range("J17").Value = 78
$ws.Range("K17").Value = -10.256410256410
$ws.Range("L17").Value = 14.754098360655
$ws.Range("M17").Value = 34.615384615384
$ws.Range("N17").Value = -9.090909090909

# Burglary (row 18)
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 43
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = -30.645161290322
$ws.Range("L18").Value = 16.216216216216
$ws.Range("M18").Value = 95.454545454545
$ws.Range("N18").Value = -66.141732283464

# Gr. Larceny (row 19)
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -43.859649122807
$ws.Range("I19").Value = 190
$ws.Range("J19").Value = 279
$ws.Range("K19").Value = -31.899641577060
$ws.Range("L19").Value = 13.772455089820
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 20.253164556962

# G.L.A. (row 20)
$ws.Range("I20").Value = 16
$ws.Range("K20").Value = -11.111111111111
$ws.Range("L20").Value = -11.111111111111
$ws.Range("M20").Value = -30.434782608695
$ws.Range("N20").Value = -88.235294117647

# TOTAL (row 21)
$ws.Range("C21").Value = 19
$ws.Range("E21").Value = -26.923076923076
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = -34.285714285714
$ws.Range("I21").Value = 368
$ws.Range("J21").Value = 503
$ws.Range("K21").Value = -26.838966202783
$ws.Range("L21").Value = 20.655737704918
$ws.Range("M21").Value = 65.022421524663
$ws.Range("N21").Value = -52.880921895006

# Transit (row 22)
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = -42.857142857142
$ws.Range("L22").Value = 33.333333333333

# Housing (row 23)
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("I23").Value = 48
$ws.Range("J23").Value = 66
$ws.Range("K23").Value = -27.272727272727
$ws.Range("L23").Value = -29.411764705882
$ws.Range("M23").Value = -12.727272727272

# Petit Larceny (row 24)
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 56
$ws.Range("E24").Value = -64.285714285714
$ws.Range("G24").Value = 174
$ws.Range("H24").Value = -52.873563218390
$ws.Range("I24").Value = 386
$ws.Range("J24").Value = 865
$ws.Range("K24").Value = -55.375722543352
$ws.Range("L24").Value = 3.208556149732
$ws.Range("M24").Value = 57.551020408163

# Misd. Assault (row 25)
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -8.333333333333
$ws.Range("I25").Value = 143
$ws.Range("J25").Value = 165
$ws.Range("K25").Value = -13.333333333333
$ws.Range("L25").Value = 57.142857142857
$ws.Range("M25").Value = 25.438596491228

# UCR Rape* (row 26) remaining plain numeric cells
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 10
$ws.Range("L26").Value = 42.857142857142

# Other Sex Crimes (row 27) remaining plain numeric cells
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 31.25
$ws.Range("L27").Value = 90.909090909090

# Hate Crimes (row 30) remaining plain numeric cell
$ws.Range("L30").Value = -72.727272727272
